$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 13 (2022, month 2): mean_PAR, max_PAR, min_PAR
$ws.Range("C13").Value = 54.7752976190476
$ws.Range("D13").Value = 662
$ws.Range("E13").Value = 0

# Update row 17 (2023, month 2): mean_PAR
$ws.Range("C17").Value = 91.8392857142857
